# Swap the presentation's theme color palette from the "Integral" scheme
# over to the stock "Office Theme" palette (the deck keeps the same 12-slot
# a:clrScheme that backs ppt/theme/theme2.xml via the slide master/Design).
#
# Target palette (Office Theme):
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

function New-RGB($r, $g, $b) {
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Order matches the clrScheme slot order exposed by Colors(1..12):
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeTheme = @(
    (New-RGB 0x00 0x00 0x00),   # dk1
    (New-RGB 0xFF 0xFF 0xFF),   # lt1
    (New-RGB 0x44 0x54 0x6A),   # dk2
    (New-RGB 0xE7 0xE6 0xE6),   # lt2
    (New-RGB 0x5B 0x9B 0xD5),   # accent1
    (New-RGB 0xED 0x7D 0x31),   # accent2
    (New-RGB 0xA5 0xA5 0xA5),   # accent3
    (New-RGB 0xFF 0xC0 0x00),   # accent4
    (New-RGB 0x44 0x72 0xC4),   # accent5
    (New-RGB 0x70 0xAD 0x47),   # accent6
    (New-RGB 0x05 0x63 0xC1),   # hlink
    (New-RGB 0x95 0x4F 0x72)    # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeTheme[$i - 1]
}
